$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was "M" stats) -> becomes "B" stats
$ws.Range("A2").Value = "B"
$ws.Range("B2").Value = 0.9459459459459459
$ws.Range("C2").Value = 0.9722222222222222
$ws.Range("D2").Value = 0.958904109589041
$ws.Range("E2").Value = 36

# Row 3 (was "B" stats) -> becomes "M" stats
$ws.Range("A3").Value = "M"
$ws.Range("B3").Value = 0.95
$ws.Range("C3").Value = 0.9047619047619048
$ws.Range("D3").Value = 0.926829268292683
$ws.Range("E3").Value = 21

# Row 4 (accuracy)
$ws.Range("B4").Value = 0.9473684210526315
$ws.Range("C4").Value = 0.9473684210526315
$ws.Range("D4").Value = 0.9473684210526315
$ws.Range("E4").Value = 0.9473684210526315

# Row 5 (macro avg)
$ws.Range("B5").Value = 0.9479729729729729
$ws.Range("C5").Value = 0.9384920634920635
$ws.Range("D5").Value = 0.942866688940862
$ws.Range("E5").Value = 57

# Row 6 (weighted avg)
$ws.Range("B6").Value = 0.9474395448079658
$ws.Range("C6").Value = 0.9473684210526315
$ws.Range("D6").Value = 0.947087062795646
$ws.Range("E6").Value = 57
